$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.536.43"
$ws.Range("E2").Value = "  +2.31%  "

$ws.Range("D3").Value = "1.987.73"
$ws.Range("E3").Value = "  +5.96%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "328.95"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "0.4688"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("D8").Value = "0.3947"

$ws.Range("D9").Value = "46.56"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "0.07968"

$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").Value = "22.81"
$ws.Range("E12").Value = "  +4.99%  "

$ws.Range("D13").Value = "2.032.78"
$ws.Range("E13").Value = "  +7.44%  "

$ws.Range("D14").Value = "7.266"
$ws.Range("E14").Value = "  +4.09%  "

$ws.Range("D15").Value = "5.886"
$ws.Range("E15").Value = "  +4.35%  "

$ws.Range("D16").Value = "0.07153"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").Value = "88.95"

$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").Value = "0.000009965"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "17.36"
$ws.Range("E20").Value = "  +2.51%  "

$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").Value = "29.644.27"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").Value = "5.555"

$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +3.22%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.121"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "158.16"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.69"
$ws.Range("E27").Value = "  +2.17%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "6.028"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "120.44"
$ws.Range("E29").Value = "  +2.85%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "1.966"
$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.09452"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.8928"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.289"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.349"
$ws.Range("E34").Value = "  +2.49%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.188"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.05845"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.178"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02133"
$ws.Range("E38").Value = "  +3.34%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "7.910"
$ws.Range("E39").Value = "  +3.18%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5764"
$ws.Range("E40").Value = "  +2.20%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.1825"
$ws.Range("E41").Value = "  +3.62%  "

$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.000003111"
$ws.Range("E42").Value = "  +97.09%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "9.824"
$ws.Range("E43").Value = "  +2.03%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "12.12"
$ws.Range("E44").Value = "  +2.66%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5385"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "2.160"
$ws.Range("E46").Value = "  -4.32%  "

$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "2.642"
$ws.Range("E47").Value = "  +5.38%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.872"
$ws.Range("E48").Value = "  +1.79%  "

$ws.Range("D49").Value = "0.06960"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "114.69"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.3125"
$ws.Range("E51").Value = "  +9.70%  "

